$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.904.42'
$ws.Range('E2').Value = '  -0.78%  '
$ws.Range('D3').Value = '2.525.46'
$ws.Range('E3').Value = '  +0.44%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '536.86'
$ws.Range('E5').Value = '  -0.79%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.17'
$ws.Range('E6').Value = '  -1.72%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.568'
$ws.Range('E8').Value = '  +0.56%  '
$ws.Range('D9').Value = '2.524.02'
$ws.Range('E9').Value = '  +0.29%  '
$ws.Range('E10').Value = '  -0.41%  '
$ws.Range('E11').Value = '  -2.25%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.31'
$ws.Range('E12').Value = '  -1.61%  '
$ws.Range('E13').Value = '  -0.88%  '
$ws.Range('D14').Value = '2.970.46'
$ws.Range('E14').Value = '  +0.21%  '
$ws.Range('E15').Value = '  -0.99%  '
$ws.Range('D16').Value = '58.925.61'
$ws.Range('E16').Value = '  -0.57%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000139'
$ws.Range('E17').Value = '  -1.26%  '
$ws.Range('D18').Value = '2.513.23'
$ws.Range('E18').Value = '  -0.02%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.17'
$ws.Range('E19').Value = '  +0.80%  '
$ws.Range('E20').Value = '  +0.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '323.91'
$ws.Range('E21').Value = '  -0.52%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('E22').Value = '  +0.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.94'
$ws.Range('E23').Value = '  +1.40%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.65'
$ws.Range('E24').Value = '  +3.31%  '
$ws.Range('E25').Value = '  -0.13%  '
$ws.Range('E26').Value = '  -2.18%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('E28').Value = '  -3.01%  '
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').Value = '0.0₃0775'
$ws.Range('E29').Value = '  -0.47%  '
$ws.Range('B30').Value = 'Aptos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.69'
$ws.Range('E30').Value = '  -0.44%  '
$ws.Range('E31').Value = '  -1.70%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '169.83'
$ws.Range('E32').Value = '  +3.75%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.17'
$ws.Range('E33').Value = '  +4.28%  '
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.44'
$ws.Range('E35').Value = '  -0.31%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.41'
$ws.Range('E36').Value = '  -0.62%  '
$ws.Range('E37').Value = '  -2.26%  '
$ws.Range('E38').Value = '  -3.12%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.72'
$ws.Range('E39').Value = '  -0.62%  '
$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.62'
$ws.Range('E40').Value = '  -1.41%  '
$ws.Range('B41').Value = 'SuiNetwork'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.811'
$ws.Range('E41').Value = '  +0.76%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '285.72'
$ws.Range('E42').Value = '  +2.35%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.11'
$ws.Range('E43').Value = '  -2.01%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '132.24'
$ws.Range('E44').Value = '  +6.12%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.997'
$ws.Range('E45').Value = '  -0.02%  '
$ws.Range('E46').Value = '  +1.89%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.87'
$ws.Range('E47').Value = '  -0.03%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0925'
$ws.Range('E48').Value = '  -1.34%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0507'
$ws.Range('E49').Value = '  -1.12%  '
$ws.Range('E50').Value = '  -1.55%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '17.36'
$ws.Range('E51').Value = '  -2.51%  '
